# daily auto push: 2026-01-21 06:54 UTC
# A new daily sample (2026/01/21, 水, 14:00, rank 201) was recorded and inserted
# into the log sheet right after the existing 2026/01/21 11:00 row (old row 695),
# which pushes every following row (old rows 696-737) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 696 - this shifts old rows 696-737 down
# to 697-738 and grows the sheet's used range to A1:D738 automatically.
$ws.Rows.Item(696).Insert()

# Fill in the newly inserted row. Column A holds dates stored as plain text
# (matching every other row in the sheet), so the value is entered with a
# leading apostrophe to stop Excel's automatic date parsing, and the style is
# then reset to "Normal" so the cell ends up as plain text with no special
# number formatting (exactly like its neighbours).
$ws.Range("A696").Value = "'2026/01/21"
$ws.Range("A696").Style = "Normal"
$ws.Range("B696").Value = "水"
$ws.Range("C696").Value = 14
$ws.Range("D696").Value = 201
